# Add 2022-Q4 data
# -----------------------------------------------------------------------
# 1. Update the "总计" (total) summary sheet: insert a new row for the
#    2022-Q4 quarter (keeping the existing 2022-Q3 / 2021-Q4 rows, which
#    shift down by one row).
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$zongji = $wb.Worksheets.Item("总计")

# Insert a new row right below the header, pushing the existing
# 2022-Q3 / 2021-Q4 rows down one position (their column-A index values
# travel with them unchanged: 0 -> row3, 1 -> row4).
$zongji.Rows.Item(2).Insert()

# The inserted row inherits the row-above formatting on every cell; only
# column A should keep that bordered/bold style, so clear it off B:D.
$zongji.Range("B2:D2").ClearFormats()

# Give cell A2 the same style used by the other index cells in column A.
$zongji.Range("A4").Copy($zongji.Range("A2"))

# Fill in the new 2022-Q4 row (first row of the table, index 0).
$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q4"
$zongji.Range("C2").Value = 5
$zongji.Range("D2").Value = 0.26

# Renumber the running index in column A for the two rows that shifted
# down: the old row 2 (index 0) is now row 3 and becomes index 1; the old
# row 3 (index 1) is now row 4 and becomes index 2.
$zongji.Range("A3").Value = 1
$zongji.Range("A4").Value = 2

# -----------------------------------------------------------------------
# 2. Create the new "2022-Q4" worksheet, positioned right before the
#    existing "2022-Q3" sheet, using that sheet's data/formatting as a
#    starting point (the two quarters share the same fund lineup), then
#    update the figures that changed quarter over quarter.
# -----------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3Index = $q3.Index
$q3.Copy($q3, $null)

# The freshly-inserted copy now occupies the position the original
# "2022-Q3" sheet used to have; grab it by that position and rename it.
$q4 = $wb.Worksheets.Item($q3Index)
$q4.Name = "2022-Q4"

# The figures below are stored as text in this workbook, so force the
# target cells to text formatting before writing the new values.
$q4.Range("D2:G6").NumberFormat = "@"

$q4.Range("D2").Value = "1.27"
$q4.Range("E2").Value = "94.90"
$q4.Range("F2").Value = "9.75"
$q4.Range("G2").Value = "0.1238"
$q4.Range("H2").Value = 3

$q4.Range("D3").Value = "0.63"
$q4.Range("E3").Value = "94.90"
$q4.Range("F3").Value = "9.75"
$q4.Range("G3").Value = "0.0614"
$q4.Range("H3").Value = 3

$q4.Range("D4").Value = "1.06"
$q4.Range("E4").Value = "82.28"
$q4.Range("F4").Value = "2.29"
$q4.Range("G4").Value = "0.0243"

$q4.Range("D5").Value = "1.06"
$q4.Range("E5").Value = "82.28"
$q4.Range("F5").Value = "2.29"
$q4.Range("G5").Value = "0.0243"

$q4.Range("D6").Value = "1.06"
$q4.Range("E6").Value = "82.28"
$q4.Range("F6").Value = "2.29"
$q4.Range("G6").Value = "0.0243"
